# Auto-generated edit script: apply scheduled-runner price updates to Titan_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J2").Value = 375
$ws.Range("L2").Value = 375
$ws.Range("N2").Value = -601
$ws.Range("H121").Value = 877.2222
$ws.Range("J121").Value = 1027.8572
$ws.Range("L121").Value = 3083.5716
$ws.Range("N121").Value = -6577.571599999999
$ws.Range("H131").Value = 9239.406999999999
$ws.Range("I131").Value = 2678.375
$ws.Range("J131").Value = 18782.727
$ws.Range("K131").Value = 8035.125
$ws.Range("L131").Value = 56348.181
$ws.Range("M131").Value = -2995.125
$ws.Range("N131").Value = -66428.181
$ws.Range("H132").Value = 21096.059
$ws.Range("I132").Value = 23058.412
$ws.Range("J132").Value = 3042.4
$ws.Range("K132").Value = 69175.236
$ws.Range("L132").Value = 9127.200000000001
$ws.Range("M132").Value = -66645.236
$ws.Range("N132").Value = -14187.2
$ws.Range("H135").Value = 978.5714
$ws.Range("I135").Value = 978.5714
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8807.142600000001
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6272.142600000001
$ws.Range("N135").ClearContents()
$ws.Range("H141").Value = 4159.4287
$ws.Range("I141").Value = 2382.3225
$ws.Range("J141").Value = 9167.637000000001
$ws.Range("K141").Value = 7146.967500000001
$ws.Range("L141").Value = 27502.911
$ws.Range("M141").Value = -1966.967500000001
$ws.Range("N141").Value = -37862.911

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35826.023
$ws.Range("I32").Value = 7817.1
$ws.Range("K32").Value = 7817.1
$ws.Range("M32").Value = -7530.1
$ws.Range("H132").Value = 3184.2307
$ws.Range("I132").Value = 2590.3103
$ws.Range("J132").Value = 4906.6
$ws.Range("K132").Value = 7770.9309
$ws.Range("L132").Value = 14719.8
$ws.Range("M132").Value = -5240.9309
$ws.Range("N132").Value = -19779.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 306.57144
$ws.Range("I22").Value = 269.4
$ws.Range("J22").Value = 399.5
$ws.Range("K22").Value = 269.4
$ws.Range("L22").Value = 399.5
$ws.Range("M22").Value = -96.39999999999998
$ws.Range("N22").Value = -745.5
$ws.Range("H105").Value = 2996.84
$ws.Range("I105").Value = 2803.889
$ws.Range("J105").Value = 3493
$ws.Range("K105").Value = 2803.889
$ws.Range("L105").Value = 3493
$ws.Range("M105").Value = -1056.889
$ws.Range("N105").Value = -6987
$ws.Range("H134").Value = 2661.3
$ws.Range("I134").Value = 2001.5834
$ws.Range("J134").Value = 4357.7144
$ws.Range("K134").Value = 6004.7502
$ws.Range("L134").Value = 13073.1432
$ws.Range("M134").Value = -3469.7502
$ws.Range("N134").Value = -18143.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 100002500
$ws.Range("J4").Value = 111113780
$ws.Range("L4").Value = 111113780
$ws.Range("N4").Value = -111114004
$ws.Range("H56").Value = 15000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H99").Value = 2899.3125
$ws.Range("I99").Value = 2127
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2127
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -629
$ws.Range("N99").Value = -6496
$ws.Range("H102").Value = 29650
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 29650
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 29650
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -34518
$ws.Range("H126").Value = 2899.3125
$ws.Range("I126").Value = 2127
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 6381
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -3911
$ws.Range("N126").Value = -15440
$ws.Range("H134").Value = 2419.889
$ws.Range("I134").Value = 1431.9688
$ws.Range("J134").Value = 4851.6924
$ws.Range("K134").Value = 4295.9064
$ws.Range("L134").Value = 14555.0772
$ws.Range("M134").Value = -1760.9064
$ws.Range("N134").Value = -19625.0772
$ws.Range("H141").Value = 283257.8
$ws.Range("J141").Value = 314366.06
$ws.Range("L141").Value = 314366.06
$ws.Range("N141").Value = -324726.06

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5651132.5
$ws.Range("J131").Value = 6804345.5
$ws.Range("L131").Value = 20413036.5
$ws.Range("N131").Value = -20423116.5
$ws.Range("H139").Value = 3817.375
$ws.Range("I139").Value = 3801.818
$ws.Range("J139").Value = 3988.5
$ws.Range("K139").Value = 11405.454
$ws.Range("L139").Value = 11965.5
$ws.Range("M139").Value = -6265.454000000002
$ws.Range("N139").Value = -22245.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4259.1304
$ws.Range("I132").Value = 3192.75
$ws.Range("J132").Value = 6696.5713
$ws.Range("K132").Value = 9578.25
$ws.Range("L132").Value = 20089.7139
$ws.Range("M132").Value = -7048.25
$ws.Range("N132").Value = -25149.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3355.2632
$ws.Range("I40").Value = 2700
$ws.Range("J40").Value = 3432.353
$ws.Range("K40").Value = 2700
$ws.Range("L40").Value = 3432.353
$ws.Range("M40").Value = -2564
$ws.Range("N40").Value = -3704.353
$ws.Range("H61").Value = 1271.9286
$ws.Range("I61").Value = 1075.25
$ws.Range("J61").Value = 1534.1666
$ws.Range("K61").Value = 1075.25
$ws.Range("L61").Value = 1534.1666
$ws.Range("M61").Value = -873.25
$ws.Range("N61").Value = -1938.1666
$ws.Range("H113").Value = 1271.9286
$ws.Range("I113").Value = 1075.25
$ws.Range("J113").Value = 1534.1666
$ws.Range("K113").Value = 1075.25
$ws.Range("L113").Value = 1534.1666
$ws.Range("M113").Value = 1094.75
$ws.Range("N113").Value = -5874.1666
$ws.Range("H122").Value = 3862.3809
$ws.Range("I122").Value = 2900
$ws.Range("K122").Value = 8700
$ws.Range("M122").Value = -6250
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 2983.4146
$ws.Range("I132").Value = 1979.5588
$ws.Range("J132").Value = 7859.2856
$ws.Range("K132").Value = 5938.6764
$ws.Range("L132").Value = 23577.8568
$ws.Range("M132").Value = -3408.6764
$ws.Range("N132").Value = -28637.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H123").Value = 33363.637
$ws.Range("J123").Value = 33363.637
$ws.Range("L123").Value = 33363.637
$ws.Range("N123").Value = -43163.637
$ws.Range("H132").Value = 3037.5
$ws.Range("I132").Value = 3258.3794
$ws.Range("J132").Value = 2610.4666
$ws.Range("K132").Value = 9775.138199999999
$ws.Range("L132").Value = 7831.399800000001
$ws.Range("M132").Value = -7245.138199999999
$ws.Range("N132").Value = -12891.3998
